# Applies the "Added many more features" edit to the Ali Baba's Gold review.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title (appears twice: Heading1 at top, bold text near the bottom)
Replace-Text "Play Ali Baba's Gold Free - Review of Top Slot Machines" "Play Ali Baba’s Gold Free Online Slot Game"

# "What we like" bullet list
Replace-Text "Stunning graphics and animation" "Breathtaking graphics and animation"
Replace-Text "Special features including respins and wild multipliers" "Special features like respins and free spins"
Replace-Text "Flexible interface allowing for perfect gameplay on any device" "Flexible interface for seamless gameplay on desktop and mobile"
Replace-Text "Rewards with high payouts" "Beautiful recreation of the Ali Baba story"

# "What we don't like" bullet list
Replace-Text "Low frequency of winnings" "High winnings but with low frequency"
Replace-Text "Limited number of paylines" "Requires patience from players"

# Closing italic summary line
Replace-Text "Read our in-depth review of Ali Baba's Gold and play for free. Enjoy stunning graphics, special features, and flexible interface for seamless gameplay." "Read our review of Ali Baba’s Gold and play this stunning online slot game for free."
